# B6-PowerPoint.pptx edit
#
# 1) The three data tables (on slides 14, 15 and 16) are switched from the
#    deck's custom "Table_0" style to the built-in "Medium Style 2 - Accent 1"
#    table style.
# 2) The design theme is switched from "Integral" (Red Violet colours) to the
#    stock "Office Theme" colour palette.

function Get-BgrFromHex([string]$hex) {
    # PowerPoint's ColorFormat.RGB takes a BGR-packed long (R + G*256 + B*65536),
    # same convention as VBA's RGB() function -- convert from a plain "RRGGBB"
    # hex string.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$newTableStyleId = "{EB572844-B146-4FF3-A7B3-9BA501FFBBA8}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the design theme's colour scheme to the stock Office palette --
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = Get-BgrFromHex $officeThemeColors[$i - 1]
}
